# Ember states and sub-states.xlsx
# - Adds the new "DoorOpen" screen-text rows (C37:C41) describing the new
#   DoorOpen sub-state screens (USB job loaded / downloading / loading /
#   error-loading / error-downloading messages), and nudges the saved
#   window position, matching the author's commit:
#   "updated states/substates spreadsheet with new DoorOpen screens, &
#    updated change log"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Best-effort: remember the workbook window position (cosmetic; not all
# hosts persist this back into bookViews, but it mirrors the source edit).
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 3120
    $win.Top = 200
} catch {
}

# C37 currently has wrap-text formatting already applied (style carried
# over from the header block) and just needs the new message text - no
# format change required.
$ws.Range("C37").Value = "<job name> loaded. Shut door to start the print."

# C38:C41 are currently blank cells using the "wrap text" variant of the
# bordered style; the new text in these rows reuses the plain (no wrap)
# bordered style already used throughout column C (e.g. C3). Copy that
# cell's format in first so the existing style is reused rather than a
# new one being created, then fill in the text.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C38:C41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C38").Value = "Downloading file..."
$ws.Range("C39").Value = "Loading file..."
$ws.Range("C40").Value = "Error loading file."
$ws.Range("C41").Value = "Error downloading file."
